$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.345.14'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').Value = '2.772.76'
$ws.Range('E3').Value = '  +1.06%  '
$ws.Range('E4').Value = '  +0.05%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '581.97'
$cell.ClearFormats()
$ws.Range('E5').Value = '  +0.62%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '161.69'
$cell.ClearFormats()
$ws.Range('E6').Value = '  +2.39%  '
$ws.Range('E7').Value = '  -0.17%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.603'
$cell.ClearFormats()
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('B10').Value = 'TRON'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '0.166'
$cell.ClearFormats()
$ws.Range('E10').Value = '  +4.82%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '5.88'
$cell.ClearFormats()
$ws.Range('E11').Value = '  +3.94%  '
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').Value = '3.263.71'
$ws.Range('E13').Value = '  +0.88%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '27.59'
$cell.ClearFormats()
$ws.Range('E14').Value = '  +2.32%  '
$ws.Range('D15').Value = '63.981.44'
$ws.Range('E15').Value = '  +0.54%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '0.0000153'
$cell.ClearFormats()
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('D17').Value = '2.778.08'
$ws.Range('E17').Value = '  +0.90%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '12.26'
$cell.ClearFormats()
$ws.Range('E18').Value = '  -0.02%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '4.88'
$cell.ClearFormats()
$ws.Range('E19').Value = '  -1.09%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '361.53'
$cell.ClearFormats()
$ws.Range('E20').Value = '  +0.48%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '6.71'
$cell.ClearFormats()
$ws.Range('E21').Value = '  -2.03%  '
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('E23').Value = '  -5.70%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '65.18'
$cell.ClearFormats()
$ws.Range('E24').Value = '  -1.56%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '0.172'
$cell.ClearFormats()
$ws.Range('E25').Value = '  +0.25%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '8.67'
$cell.ClearFormats()
$ws.Range('E26').Value = '  +0.23%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.ClearFormats()
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('D28').Value = '0.0₃0923'
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('E29').Value = '  +4.97%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '2.00'
$cell.ClearFormats()
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '1.38'
$cell.ClearFormats()
$ws.Range('E31').Value = '  +11.54%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '167.00'
$cell.ClearFormats()
$ws.Range('E32').Value = '  -1.54%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '1.54'
$cell.ClearFormats()
$ws.Range('E33').Value = '  +5.00%  '
$ws.Range('E34').Value = '  +0.64%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '20.28'
$cell.ClearFormats()
$ws.Range('E35').Value = '  -1.16%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.ClearFormats()
$ws.Range('E36').Value = '  -0.01%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '1.84'
$cell.ClearFormats()
$ws.Range('E37').Value = '  +1.83%  '
$ws.Range('E38').Value = '  +1.42%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '352.53'
$cell.ClearFormats()
$ws.Range('E39').Value = '  +6.22%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '6.42'
$cell.ClearFormats()
$ws.Range('E40').Value = '  +5.28%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '4.22'
$cell.ClearFormats()
$ws.Range('E41').Value = '  +0.82%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '39.35'
$cell.ClearFormats()
$ws.Range('E42').Value = '  -0.79%  '
$ws.Range('E43').Value = '  +4.32%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '21.72'
$cell.ClearFormats()
$ws.Range('E44').Value = '  -0.89%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '0.0598'
$cell.ClearFormats()
$ws.Range('E45').Value = '  +0.39%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '137.96'
$cell.ClearFormats()
$ws.Range('E46').Value = '  +1.22%  '
$ws.Range('E47').Value = '  -0.85%  '
$ws.Range('E48').Value = '  -0.77%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '0.102'
$cell.ClearFormats()
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('E50').Value = '  -0.34%  '
$ws.Range('D51').Value = '2.146.85'
$ws.Range('E51').Value = '  +1.55%  '
